$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.502.28'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '2.464.31'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.13'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.25'
$ws.Range("E6").Value = '  -4.00%  '
$ws.Range("E7").Value = '  -3.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.486'
$ws.Range("E9").Value = '  -4.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '31.91'
$ws.Range("E10").Value = '  -6.08%  '
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("D13").Value = '2.841.96'
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("E14").Value = '  -3.77%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.454.48'
$ws.Range("E15").Value = '  -2.92%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.10'
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.756'
$ws.Range("E17").Value = '  -4.29%  '
$ws.Range("D18").Value = '41.338.29'
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.21'
$ws.Range("E19").Value = '  -2.79%  '
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.87'
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.89'
$ws.Range("E22").Value = '  -6.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.06'
$ws.Range("E23").Value = '  -1.84%  '
$ws.Range("E24").Value = '  -4.05%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.86'
$ws.Range("E26").Value = '  -3.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.96'
$ws.Range("E27").Value = '  -3.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.23'
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.50'
$ws.Range("E29").Value = '  -2.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.66'
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '151.53'
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.30'
$ws.Range("E32").Value = '  -5.80%  '
$ws.Range("E33").Value = '  -3.47%  '
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.67'
$ws.Range("E35").Value = '  +2.37%  '
$ws.Range("E36").Value = '  -3.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.93'
$ws.Range("E37").Value = '  -3.45%  '
$ws.Range("E38").Value = '  -4.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.112'
$ws.Range("E39").Value = '  -2.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0999'
$ws.Range("E40").Value = '  -6.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.03'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.21'
$ws.Range("E43").Value = '  -4.74%  '
$ws.Range("D44").Value = '1.939.23'
$ws.Range("E44").Value = '  -2.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0277'
$ws.Range("E45").Value = '  -3.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.91'
$ws.Range("E46").Value = '  -5.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.67'
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("D48").Value = '2.705.80'
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '94.39'
$ws.Range("E49").Value = '  -3.56%  '
$ws.Range("E50").Value = '  -4.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.58'
$ws.Range("E51").Value = '  -5.07%  '
